{"js": "// Update the author byline on the Synopsis page: add \"Matthew N. White\"\n// to the existing author list (\"Christopher Carroll, Alan Lujan\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst oldText = \"Christopher Carroll, Alan Lujan\";\nconst newText = \"Christopher Carroll, Alan Lujan, Matthew N. White\";\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.style === \"Heading 4\" && p.text === oldText) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  // Fallback: search the whole document body for the exact author text.\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n} else {\n  target.insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the author byline on the Synopsis page: add \"Matthew N. White\"\n# to the existing author list (\"Christopher Carroll, Alan Lujan\").\n$doc = $word.ActiveDocument\n\n$oldText = \"Christopher Carroll, Alan Lujan\"\n$newText = \"Christopher Carroll, Alan Lujan, Matthew N. White\"\n\n$find = $doc.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
